$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '41.984.04'
$ws.Cells.Item(2, 5).Value = '  -0.62%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.211.06'
$ws.Cells.Item(3, 5).Value = '  -1.50%  '

$ws.Cells.Item(4, 5).Value = '  +0.07%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '240.53'
$ws.Cells.Item(5, 5).Value = '  -2.67%  '

$ws.Cells.Item(6, 5).Value = '  -1.06%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '72.88'
$ws.Cells.Item(7, 5).Value = '  -2.87%  '

$ws.Cells.Item(8, 5).Value = '  +0.10%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.601'
$ws.Cells.Item(9, 5).Value = '  -2.99%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '42.23'
$ws.Cells.Item(10, 5).Value = '  -0.21%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0947'
$ws.Cells.Item(11, 5).Value = '  +0.13%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.02'
$ws.Cells.Item(12, 5).Value = '  -2.10%  '

$ws.Cells.Item(13, 5).Value = '  -0.69%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '2.543.86'
$ws.Cells.Item(14, 5).Value = '  -1.25%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '14.16'
$ws.Cells.Item(15, 5).Value = '  -2.57%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.830'
$ws.Cells.Item(16, 5).Value = '  -2.86%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.222.33'
$ws.Cells.Item(17, 5).Value = '  +0.13%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '41.822.23'
$ws.Cells.Item(18, 5).Value = '  -0.74%  '

$ws.Cells.Item(19, 5).Value = '  +8.61%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '72.79'
$ws.Cells.Item(20, 5).Value = '  +0.78%  '

$ws.Cells.Item(21, 5).Value = '  -1.10%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '10.30'
$ws.Cells.Item(22, 5).Value = '  +15.40%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '228.42'
$ws.Cells.Item(23, 5).Value = '  -1.58%  '

$ws.Cells.Item(24, 5).Value = '  -4.79%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '11.58'
$ws.Cells.Item(25, 5).Value = '  +2.24%  '

$ws.Cells.Item(27, 5).Value = '  -0.78%  '

$ws.Cells.Item(28, 5).Value = '  -2.05%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.17'
$ws.Cells.Item(29, 5).Value = '  +0.86%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '166.65'
$ws.Cells.Item(30, 5).Value = '  -1.56%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '20.49'
$ws.Cells.Item(31, 5).Value = '  -0.91%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '5.71'
$ws.Cells.Item(32, 5).Value = '  +9.56%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0787'
$ws.Cells.Item(33, 5).Value = '  -4.21%  '

$ws.Cells.Item(34, 5).Value = '  -0.99%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '28.72'
$ws.Cells.Item(35, 5).Value = '  -6.90%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.110'
$ws.Cells.Item(36, 5).Value = '  -8.54%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.23'
$ws.Cells.Item(37, 5).Value = '  -5.81%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0298'
$ws.Cells.Item(38, 5).Value = '  -4.48%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '13.31'
$ws.Cells.Item(39, 5).Value = '  -3.30%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '65.43'
$ws.Cells.Item(40, 5).Value = '  +5.02%  '

$ws.Cells.Item(41, 5).Value = '  -4.03%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.60'
$ws.Cells.Item(42, 5).Value = '  -2.98%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.197'
$ws.Cells.Item(43, 5).Value = '  -3.90%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '8.62'
$ws.Cells.Item(44, 5).Value = '  -0.86%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '103.48'
$ws.Cells.Item(45, 5).Value = '  -1.66%  '

$ws.Cells.Item(46, 5).Value = '  -2.32%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.36'
$ws.Cells.Item(47, 5).Value = '  +3.20%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.11'
$ws.Cells.Item(48, 5).Value = '  -1.46%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.15'
$ws.Cells.Item(49, 5).Value = '  -0.91%  '

$ws.Cells.Item(50, 5).Value = '  -0.05%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.417.13'
$ws.Cells.Item(51, 5).Value = '  -2.24%  '
